# Apply the "Added periods to exclude from GAGE" edit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. CGO: add new excluded-period rows (CH4 + CFC-113), with comments
# ---------------------------------------------------------------------
$cgo = $wb.Worksheets.Item("CGO")

$cgo.Range("A5").Value = "ch4"
$cgo.Range("B5").Value = "GAGE"
$cgo.Range("C5").Value = "1990-01-08 12:09"
$cgo.Range("D5").Value = "1990-01-17 06:04"

$cgo.Range("A6").Value = "ch4"
$cgo.Range("B6").Value = "GAGE"
$cgo.Range("C6").Value = "1990-02-02 21:57"
$cgo.Range("D6").Value = "1990-02-05 06:53"

$cgo.Range("A7").Value = "cfc-113"
$cgo.Range("B7").Value = "GAGE"
$cgo.Range("C7").Value = "1982-09-05 23:10"
$cgo.Range("D7").Value = "1982-09-05 23:20"

$cgo.Range("A8").Value = "cfc-113"
$cgo.Range("B8").Value = "GAGE"
$cgo.Range("C8").Value = "1982-10-16 14:40"
$cgo.Range("D8").Value = "1982-10-16 14:50"

$cgo.Range("A9").Value = "cfc-113"
$cgo.Range("B9").Value = "GAGE"
$cgo.Range("C9").Value = "1982-10-25 14:25"
$cgo.Range("D9").Value = "1982-10-25 14:35"

$cgo.Range("A10").Value = "cfc-113"
$cgo.Range("B10").Value = "GAGE"
$cgo.Range("C10").Value = "1982-12-28 17:45"
$cgo.Range("D10").Value = "1982-12-28 18:55"

$ch4Comment = "Ray Wang (pers comm., 31/7/23): `n1), CH4 data from 1990/01/08, 23:09 to 1990/01/17, 17:04 (YYYY/mm/dd, hh:mm).   `n2), CH4 data from 1990/02/03, 08:57 to 1990/02/05, 17:53.`nCONVERTED TO UTC FROM LOCAL TIME"
$cgo.Range("C5").AddCommentThreaded($ch4Comment) | Out-Null

$cfcComment = "Paul Krummel (pers comm., 31/7/23): point to flag please:`nCGO GAGE CFC-113 1982 09 06 0916 (local time) or 1982 09 05 2316 (UTC time)"
$cgo.Range("C7").AddCommentThreaded($cfcComment) | Out-Null

# ---------------------------------------------------------------------
# 2. New worksheets: RPB, CMO, SMO (added after GSN, in that order)
# ---------------------------------------------------------------------
$gsn = $wb.Worksheets.Item("GSN")

$rpb = $wb.Worksheets.Add($null, $gsn)
$rpb.Name = "RPB"

$cmo = $wb.Worksheets.Add($null, $rpb)
$cmo.Name = "CMO"

$smo = $wb.Worksheets.Add($null, $cmo)
$smo.Name = "SMO"

foreach ($ws in @($rpb, $cmo, $smo)) {
    $ws.Range("A1").NumberFormat = "@"
    $ws.Range("A1").Value = "# Mole fraction data during the ranges in this worksheet will be excluded"
    $ws.Range("A2").NumberFormat = "@"
    $ws.Range("A2").Value = "# Date format must by YYYY-MM-DD HH:MM"
    $ws.Range("A3").NumberFormat = "@"
    $ws.Range("A3").Value = "# Make sure that the cell format is text rather than Excel's date format"
    $ws.Range("A4:D4").NumberFormat = "@"
    $ws.Range("A4").Value = "Species"
    $ws.Range("B4").Value = "Instrument"
    $ws.Range("C4").Value = "Start"
    $ws.Range("D4").Value = "End"
    $ws.Range("A5:D5").NumberFormat = "@"
}

# --- RPB: remove all GAGE CH4 data ---
$rpb.Range("A5").Value = "ch4"
$rpb.Range("B5").Value = "GAGE"
$rpb.Range("C5").Value = "1970-01-01 00:00"
$rpb.Range("D5").Value = "2050-01-01 00:00"
$rpb.Range("C5").AddCommentThreaded("Ray Wang (pers comm., 31/7/23): Remove all GAGE CH4 data from RPB") | Out-Null

# --- CMO (Oregon): remove CH4 data for Sept/Oct 1985 ---
$cmo.Range("A5").Value = "ch4"
$cmo.Range("B5").Value = "GAGE"
$cmo.Range("C5").Value = "1985-09-01 00:00"
$cmo.Range("D5").Value = "1985-11-01 00:00"
$cmo.Range("C5").AddCommentThreaded("Ray Wang (pers comm., 31/7/23): Remove all CH4 data in sept and oct 1985 from Oregon") | Out-Null

# --- SMO: remove N2O after Dec 1995 ---
$smo.Range("A5").Value = "n2o"
$smo.Range("B5").Value = "GAGE"
$smo.Range("C5").Value = "1995-12-31 00:00"
$smo.Range("D5").Value = "2050-01-01 00:00"
$smo.Range("C5").AddCommentThreaded("Ray Wang (pers comm., 31/7/23): remove all N2O after Dec. 1995 (i.e. from 1996/01/01 to the end of GAGE period).") | Out-Null

# ---------------------------------------------------------------------
# 3. Selection / active-tab bookkeeping to match the saved view state
# ---------------------------------------------------------------------
$rpb.Range("E18").Select() | Out-Null
$cmo.Range("A6").Select() | Out-Null
$smo.Range("E21").Select() | Out-Null

$cgo.Activate() | Out-Null
$cgo.Range("E18").Select() | Out-Null
